$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARN: not found -> $find"
    }
}

# 1. Title: "Initial research" -> "Research Plan"
ReplaceText "Initial research" "Research Plan"

# 2. Date day: 04 -> 09
ReplaceText "Date: 04-03-2024" "Date: 09-03-2024"

# 3. Introduction purpose paragraph
ReplaceText "The purpose of this document is to outline the research process conducted in order to determine what…." "The purpose of this document is to outline the planning for the research that I am going to undertake during the creation of the HeardIT application. In this document I will outline the main research question and the following sub-questions that will be established in order to answer the main question."

# 4. Main research question
ReplaceText "How to improve and further extend the testing infrastructure around the current testing environment for the Loopsorter stream?" "How to ensure that HeardIT meets the modern standards for enterprise level software applications?"

# 5. "Answering this question..." paragraph
ReplaceText "Answering this question will allow us to find what are some of the main issues, that are within the scope of my assignment, that the Loopsorter stream is facing with their testing environment." "Answering this question will allow us to definitively prove that HeardIT is an application that meets the modern demands and standards for software applications. To answer this question, sub-questions that focus on specific aspects from the main question are established. Finding the answer to those will allow me to conclude a definitive answer for the main question. In the next section, I will establish the sub-questions and provide an answer to each of them."

# 6. Sub-question 1
ReplaceText "How can test engineers determine what kind of failures occur for the Test-Broker Engine and how often?" "What is the best architecture type that will allow me to create a web-application that meets the modern standards? "

# 7. Sub-question 2
ReplaceText "How can test engineers find the root cause of a failed build more easily?" "What are the best technologies that I should use when developing HeardIT?"

# 8. Sub-question 3
ReplaceText "What are some of the most common issues that occur?" "How can I test to ensure that HeardIT meets the standards of modern applications?"

# 9. "At this point ..." paragraph trimmed
ReplaceText "At this point I already had created the main part of " "At this point "

Write-Output "done-phase1"

# 10. Move the "Methods used: Literature study, Expert interview, Interview" paragraph
#     so that it appears right after the empty paragraph following sub-question 1
#     (i.e. before "In order to answer this question, I firstly "), leaving the old
#     location as an empty paragraph (same pPr, no runs).
$targetIdx = -1
$sourceIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "In order to answer this question*") { $targetIdx = $i }
    if ($t -like "Methods used: Literature study*") { $sourceIdx = $i }
}
Write-Output "targetIdx=$targetIdx sourceIdx=$sourceIdx"
if ($targetIdx -gt 0 -and $sourceIdx -gt 0) {
    $ptarget = $d.Paragraphs($targetIdx)
    $ptarget.Range.InsertParagraphBefore()
    $newp = $d.Paragraphs($targetIdx)
    $newp.Format.LeftIndent = 18
    $newp.Format.FirstLineIndent = 18
    $newp.Range.Text = "Methods used: Literature study, Expert interview, Interview"

    # the source paragraph shifted down by 1 because of the insertion
    $oldp = $d.Paragraphs($sourceIdx + 1)
    $r = $oldp.Range
    $r.MoveEnd(1, -1)
    $r.Text = ""
}

Write-Output "done-phase2"
